# Insert two new data rows right before the current row 1074.
# This mirrors a weekly data refresh where two new observations were
# prepended to this block of records, pushing all of the existing rows
# at 1074-1155 down to 1076-1157 (dimension grows from R1155 to R1157).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 1074 (one at a time, since inserting
# a combined multi-row range does not behave as a true row-shift here).
$ws.Rows.Item(1074).Insert()
$ws.Rows.Item(1074).Insert()

# First new row (becomes row 1074) - "Primera" quality observation.
$ws.Cells.Item(1074, 1).Value2 = 6
$ws.Cells.Item(1074, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1074, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1074, 4).Value2 = 45013
$ws.Cells.Item(1074, 5).Value2 = 13
$ws.Cells.Item(1074, 6).Value2 = 100112009
$ws.Cells.Item(1074, 7).Value2 = "Acelga"
$ws.Cells.Item(1074, 8).Value2 = "Sin especificar"
$ws.Cells.Item(1074, 9).Value2 = "Primera"
$ws.Cells.Item(1074, 10).Value2 = 360
$ws.Cells.Item(1074, 11).Value2 = 12000
$ws.Cells.Item(1074, 12).Value2 = 13000
$ws.Cells.Item(1074, 13).Value2 = 12472
$ws.Cells.Item(1074, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(1074, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(1074, 16).Value2 = 4157
$ws.Cells.Item(1074, 17).Value2 = 3
$ws.Cells.Item(1074, 18).Value2 = "Hortaliza"

# Second new row (becomes row 1075) - "Segunda" quality observation.
$ws.Cells.Item(1075, 1).Value2 = 6
$ws.Cells.Item(1075, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1075, 3).Value2 = "Metropolitana"
$ws.Cells.Item(1075, 4).Value2 = 45013
$ws.Cells.Item(1075, 5).Value2 = 13
$ws.Cells.Item(1075, 6).Value2 = 100112009
$ws.Cells.Item(1075, 7).Value2 = "Acelga"
$ws.Cells.Item(1075, 8).Value2 = "Sin especificar"
$ws.Cells.Item(1075, 9).Value2 = "Segunda"
$ws.Cells.Item(1075, 10).Value2 = 120
$ws.Cells.Item(1075, 11).Value2 = 9000
$ws.Cells.Item(1075, 12).Value2 = 9000
$ws.Cells.Item(1075, 13).Value2 = 9000
$ws.Cells.Item(1075, 14).Value2 = "$/docena de atados"
$ws.Cells.Item(1075, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(1075, 16).Value2 = 3000
$ws.Cells.Item(1075, 17).Value2 = 3
$ws.Cells.Item(1075, 18).Value2 = "Hortaliza"
